# ---------------------------------------------------------------------------
# "trying to profile memory and stuff"
#
# - Updates the measured timing numbers on the two existing raw-timings
#   sheets ("Data Set 0 Timings (Pd)" and "Data Set 0 Timings (TD)").
# - Inserts a new "regression - DS 1, T_n 16 (Pd)" sheet right after the Pd
#   timings sheet, and a new "regression - DS 1, T_n 16 (TD)" sheet right
#   after the TD timings sheet, each holding a small regression-coefficient
#   table (raw data MRL / subset MRL per predictor).
# - Refreshes the "Data Set 0 Timings (combined)" sheet (Pd rows followed by
#   TD rows) with the same updated numbers.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsPd   = $wb.Worksheets.Item(1)
$wsTd   = $wb.Worksheets.Item(2)
$wsComb = $wb.Worksheets.Item(3)

# New measured values for the Pd timings sheet, columns F,G,H,I,J,L (rows 2-6).
$pdRows = @(
    @(11.468, 0.314,  7.771, 1.65,  1.653, 0.00082),
    @(12.251, 0.329,  8.461, 1.652, 1.685, 0.00086),
    @(12.68,  0.356,  8.779999999999999, 1.646, 1.81, 0.00101),
    @(14.29,  0.424, 10.034, 1.855, 1.837, 0.0009700000000000001),
    @(13.805, 0.393,  9.692, 1.782, 1.839, 0.0009300000000000001)
)

# New measured values for the TD timings sheet, columns F,G,H,I,J,L (rows 2-6).
$tdRows = @(
    @(111.256, 1.017, 1.163, 46.593, 46.294, 0.00095),
    @(111.614, 1.013, 1.181, 46.402, 47.963, 0.00098),
    @(114.423, 1.022, 1.194, 47.634, 47.448, 0.0009700000000000001),
    @(112.195, 1.016, 1.166, 47.836, 47.101, 0.00098),
    @(123.969, 1.032, 1.192, 50.912, 55.727, 0.0009700000000000001)
)

function Set-TimingRows($ws, $rows, $startRow) {
    for ($i = 0; $i -lt $rows.Count; $i++) {
        $r = $startRow + $i
        $vals = $rows[$i]
        $ws.Cells.Item($r, 6).Value  = $vals[0]   # F - Setup Time
        $ws.Cells.Item($r, 7).Value  = $vals[1]   # G - c1_t
        $ws.Cells.Item($r, 8).Value  = $vals[2]   # H - c2_t
        $ws.Cells.Item($r, 9).Value  = $vals[3]   # I - c3_t
        $ws.Cells.Item($r, 10).Value = $vals[4]   # J - c4_t
        $ws.Cells.Item($r, 12).Value = $vals[5]   # L - Solve Time
    }
}

# Update the two raw-timings sheets in place.
Set-TimingRows $wsPd   $pdRows 2
Set-TimingRows $wsTd   $tdRows 2

# Update the combined sheet: Pd rows (2-6) then TD rows (7-11), same numbers.
Set-TimingRows $wsComb $pdRows 2
Set-TimingRows $wsComb $tdRows 7

# ---------------------------------------------------------------------------
# Regression-coefficient tables (identical content for both new sheets).
# Column B = "raw data MRL", Column C = "subset MRL".
# ---------------------------------------------------------------------------
$regRows = @(
    @("IMPLANT",               -1.645274458028569,   -23.21707356636094),
    @("B_MOP",                  -3.919640654934591,   -2.553320925882868),
    @("B_COP",                  -3.69113390607585,     0.132416941088189),
    @("B_COC",                  10.79838647173405,     [double]"-2.831068712794149e-15"),
    @("B_DM",                   -3.187611910723564,     2.420903984794039),
    @("POLY_UHWMPE",           -15.64469938464845,     -2.189579731453116),
    @("POLY_XPLE",               2.273419237814216,     2.382537352323087),
    @("POLY_A_XPLE",             2.572893675100173,    -0.1929576208699801),
    @("HEAD_22mm",               [double]"1.77635683940025e-15", [double]"1.110223024625157e-16"),
    @("HEAD_28mm",              -7.290197463123523,     5.243063723572009),
    @("HEAD_32mm",               2.605843473573495,    -2.694321738188861),
    @("HEAD_36mm",               3.805859788548451,     0.3998858734049386),
    @("HEAD_40mm",               5.755364303716566,     1.663714769400042),
    @("HEAD_44mm",              -4.876870102714975,    -4.612342628188121),
    @("APP_anterior",            1.881853646893045,     1.686073766131556),
    @("APP_anterolateral",       3.518868122288081,     2.836849038998879),
    @("APP_posterior",           3.009654668471327,    -2.431930382566656),
    @("APP_transtrochanteric",  -8.410376437652463,    -2.09099242256379),
    @("S_VOLLUME",              -0.02701493435443902,  -0.5972006935788479),
    @("FEMALE",                  0.6301467229078077,    1.055938687732325),
    @("BMI",                     0.08245388996230396,  -0.1215738358437359),
    @("const",                  59.55094637106095,    112.8474444374117),
    @("R^2 Score",               0.05408633511067074,   0.4413894059322662)
)

function Fill-RegressionSheet($ws) {
    $ws.Range("B1").Value = "raw data MRL"
    $ws.Range("C1").Value = "subset MRL"

    for ($i = 0; $i -lt $regRows.Count; $i++) {
        $r = $i + 2
        $row = $regRows[$i]
        $ws.Cells.Item($r, 1).Value = $row[0]
        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 3).Value = $row[2]
    }

    # Match the bold/centered/bordered style used for headers & row labels
    # elsewhere in this workbook (copy format only, values already set above).
    $wsPd.Range("B1:C1").Copy()
    $ws.Range("B1:C1").PasteSpecial(-4122)
    $wsPd.Range("A2").Copy()
    $ws.Range("A2:A24").PasteSpecial(-4122)
}

# Insert the two new sheets. Do this back-to-front (TD-side first) so that
# earlier indices/object references used below stay valid as sheets shift.
$regTd = $wb.Worksheets.Add($null, $wsTd)
$regTd.Name = "regression - DS 1, T_n 16 (TD)"
Fill-RegressionSheet $regTd

$regPd = $wb.Worksheets.Add($null, $wsPd)
$regPd.Name = "regression - DS 1, T_n 16 (Pd)"
Fill-RegressionSheet $regPd

$wsPd.Select()
